$wb = $excel.ActiveWorkbook

# --- Replace "Sheet2" with a fresh "Sheet1" ---------------------------------
# Adding a worksheet when "Sheet2" already exists makes Excel name the new,
# blank sheet "Sheet1" and place it first; deleting the old "Sheet2"
# afterwards leaves a single sheet named "Sheet1" whose sheetId has been
# bumped (1 -> 2), matching <sheet name="Sheet1" sheetId="2" .../>.
$new = $wb.Worksheets.Add()
$old = $wb.Worksheets.Item(2)
$old.Delete() | Out-Null

$ws = $wb.ActiveSheet

# --- Header row (row 3, columns C..K) --------------------------------------
$ws.Range("C3").Value = "Name"
$ws.Range("D3").Value = "Emailid"
$ws.Range("E3").Value = "ContactNo"
$ws.Range("F3").Value = "Designation"
$ws.Range("G3").Value = "Address"
$ws.Range("H3").Value = "City"
$ws.Range("I3").Value = "State"
$ws.Range("J3").Value = "Country"
$ws.Range("K3").Value = "Dateofbirth"
$ws.Range("C3:K3").Font.Bold = $true

# --- Data / placeholder row (row 4) -----------------------------------------
$ws.Range("C4").Value = "XXXXXX"

# Email placeholder with a mailto hyperlink
$ws.Range("D4").Value = "xxxxx@xxxx.xxx"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:xxxxx@xxxx.xxx") | Out-Null

# Date of birth placeholder cell - date-formatted, left blank
$ws.Range("K4").NumberFormat = "mm-dd-yy"

# --- Column widths -----------------------------------------------------------
$ws.Range("C1").ColumnWidth = 20.28515625
$ws.Range("D1").ColumnWidth = 32.7109375
$ws.Range("E1").ColumnWidth = 25.42578125
$ws.Range("F1").ColumnWidth = 15.42578125
$ws.Range("G1").ColumnWidth = 19.140625
$ws.Range("K1").ColumnWidth = 14.5703125

# --- Selection matches the saved view ----------------------------------------
$ws.Range("C3:K4").Select()
